{"js": "// Applies the \"Physics todo\" list edits described by the target diff:\n//  1. \"Physics todo:\" -> re-split into runs with a proofed (\"todo\") spell-check span (text unchanged).\n//  2. \"Add polygon shape collision\" -> append \" (GJK)\"; add a new empty bullet paragraph after it.\n//  3. Insert a new bullet \"Add button to pause and iterate through\" right before\n//     \"Check if inversetransformpoint actually works\", and mark \"inversetransformpoint\" as proofed.\n//  4. \"Fix the polygon aabb, it aint global rn\" and \"Make ui prettier \" -> re-split into runs with\n//     proofed spans around \"aabb\", \"aint\", \"rn\", \"ui\" (text unchanged).\n//  5. Delete the bullet \"Uncapping framerate for speedy physics\".\n\nconst WORD_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\n// Wrap one or more <w:p> elements (as a raw string) in the minimal OOXML package\n// that Range.insertOoxml expects.\nfunction pkg(bodyInnerXml) {\n  return `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document ${WORD_NS}>\n        <w:body>\n          ${bodyInnerXml}\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n}\n\nconst LIST_PPR = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>';\n\n// Find a (non-deleted) paragraph anchored on its exact current text.\nasync function findParagraph(body, exactText) {\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n  for (const p of paras.items) {\n    if (p.text === exactText) return p;\n  }\n  throw new Error(\"paragraph not found: \" + JSON.stringify(exactText));\n}\n\n// 1) \"Physics todo:\" -> split into 3 runs, proofing \"todo\".\n{\n  const p = await findParagraph(context.document.body, \"Physics todo:\");\n  const range = p.getRange();\n  const xml = pkg(`<w:p>\n    <w:r><w:t xml:space=\"preserve\">Physics </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>todo</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t>:</w:t></w:r>\n  </w:p>`);\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"Add polygon shape collision\" -> append \" (GJK)\"; add empty bullet paragraph after it.\n{\n  const p = await findParagraph(context.document.body, \"Add polygon shape collision\");\n  const range = p.getRange();\n  const xml = pkg(`<w:p>\n    ${LIST_PPR}\n    <w:r><w:t>Add polygon shape collision</w:t></w:r>\n    <w:r><w:t xml:space=\"preserve\"> (GJK)</w:t></w:r>\n  </w:p>\n  <w:p>\n    ${LIST_PPR}\n  </w:p>`);\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Insert \"Add button to pause and iterate through\" before the \"Check if ...\" bullet,\n//    and mark \"inversetransformpoint\" as proofed within that bullet.\n{\n  const p = await findParagraph(context.document.body, \"Check if inversetransformpoint actually works\");\n  const range = p.getRange();\n  const xml = pkg(`<w:p>\n    ${LIST_PPR}\n    <w:r><w:t>Add button to pause and iterate through</w:t></w:r>\n  </w:p>\n  <w:p>\n    ${LIST_PPR}\n    <w:r><w:t xml:space=\"preserve\">Check if </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>inversetransformpoint</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\"> actually works</w:t></w:r>\n  </w:p>`);\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 4) \"Fix the polygon aabb, it aint global rn\" -> proof \"aabb\", \"aint\", \"rn\".\n{\n  const p = await findParagraph(context.document.body, \"Fix the polygon aabb, it aint global rn\");\n  const range = p.getRange();\n  const xml = pkg(`<w:p>\n    ${LIST_PPR}\n    <w:r><w:t xml:space=\"preserve\">Fix the polygon </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>aabb</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\">, it </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>aint</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\"> global</w:t></w:r>\n    <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>rn</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n  </w:p>`);\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \"Make ui prettier \" -> proof \"ui\".\n{\n  const p = await findParagraph(context.document.body, \"Make ui prettier \");\n  const range = p.getRange();\n  const xml = pkg(`<w:p>\n    ${LIST_PPR}\n    <w:r><w:t xml:space=\"preserve\">Make </w:t></w:r>\n    <w:proofErr w:type=\"spellStart\"/>\n    <w:r><w:t>ui</w:t></w:r>\n    <w:proofErr w:type=\"spellEnd\"/>\n    <w:r><w:t xml:space=\"preserve\"> prettier</w:t></w:r>\n    <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  </w:p>`);\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 5) Delete the bullet \"Uncapping framerate for speedy physics\".\n{\n  const p = await findParagraph(context.document.body, \"Uncapping framerate for speedy physics\");\n  p.delete();\n  await context.sync();\n}\n", "ps1": "# Applies the \"Physics todo\" list edits described by the target diff:\n#  1. \"Physics todo:\" -> re-split into runs with a proofed (\"todo\") spell-check span (text unchanged).\n#  2. \"Add polygon shape collision\" -> append \" (GJK)\"; add a new empty bullet paragraph after it.\n#  3. Insert a new bullet \"Add button to pause and iterate through\" right before\n#     \"Check if inversetransformpoint actually works\", and mark \"inversetransformpoint\" as proofed.\n#  4. \"Fix the polygon aabb, it aint global rn\" and \"Make ui prettier \" -> re-split into runs with\n#     proofed spans around \"aabb\", \"aint\", \"rn\", \"ui\" (text unchanged).\n#  5. Delete the bullet \"Uncapping framerate for speedy physics\".\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphByText($doc, [string]$text) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -eq ($text + \"`r\")) {\n            return $p\n        }\n    }\n    throw \"paragraph not found: $text\"\n}\n\n$listPPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>'\n\nfunction New-PkgXml([string]$bodyInnerXml) {\n    return @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          $bodyInnerXml\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n}\n\n# 1) \"Physics todo:\" -> split into 3 runs, proofing \"todo\".\n$p = Find-ParagraphByText $d \"Physics todo:\"\n$xml = New-PkgXml @\"\n<w:p>\n  <w:r><w:t xml:space=\"preserve\">Physics </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>todo</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t>:</w:t></w:r>\n</w:p>\n\"@\n$null = $p.Range.InsertXML($xml)\n\n# 2) \"Add polygon shape collision\" -> append \" (GJK)\"; add empty bullet paragraph after it.\n$p = Find-ParagraphByText $d \"Add polygon shape collision\"\n$xml = New-PkgXml @\"\n<w:p>\n  $listPPr\n  <w:r><w:t>Add polygon shape collision</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> (GJK)</w:t></w:r>\n</w:p>\n<w:p>\n  $listPPr\n</w:p>\n\"@\n$null = $p.Range.InsertXML($xml)\n\n# 3) Insert \"Add button to pause and iterate through\" before the \"Check if ...\" bullet,\n#    and mark \"inversetransformpoint\" as proofed within that bullet.\n$p = Find-ParagraphByText $d \"Check if inversetransformpoint actually works\"\n$xml = New-PkgXml @\"\n<w:p>\n  $listPPr\n  <w:r><w:t>Add button to pause and iterate through</w:t></w:r>\n</w:p>\n<w:p>\n  $listPPr\n  <w:r><w:t xml:space=\"preserve\">Check if </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>inversetransformpoint</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> actually works</w:t></w:r>\n</w:p>\n\"@\n$null = $p.Range.InsertXML($xml)\n\n# 4) \"Fix the polygon aabb, it aint global rn\" -> proof \"aabb\", \"aint\", \"rn\".\n$p = Find-ParagraphByText $d \"Fix the polygon aabb, it aint global rn\"\n$xml = New-PkgXml @\"\n<w:p>\n  $listPPr\n  <w:r><w:t xml:space=\"preserve\">Fix the polygon </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>aabb</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\">, it </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>aint</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> global</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>rn</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n</w:p>\n\"@\n$null = $p.Range.InsertXML($xml)\n\n# \"Make ui prettier \" -> proof \"ui\".\n$p = Find-ParagraphByText $d \"Make ui prettier \"\n$xml = New-PkgXml @\"\n<w:p>\n  $listPPr\n  <w:r><w:t xml:space=\"preserve\">Make </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>ui</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> prettier</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n</w:p>\n\"@\n$null = $p.Range.InsertXML($xml)\n\n# 5) Delete the bullet \"Uncapping framerate for speedy physics\".\n$p = Find-ParagraphByText $d \"Uncapping framerate for speedy physics\"\n$null = $p.Range.Delete()\n"}
